# Generate Report for Handback
# - Update the "Status" text for the c2276190-... row from
#   "Ready for handoff" to "Handback transform failed" (shared across the
#   Overview, zh-cn and de-de sheets since it's a shared string).
# - Populate the "Error Detail" column (P) for that same row on the
#   zh-cn and de-de sheets with a handback-mismatch diagnostic message.
# - Widen the "Error Detail" column on those two sheets to fit the text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# The Status column value ("Ready for handoff" -> "Handback transform
# failed") shows up in four cells: Overview!E3, Overview!F3, zh-cn!C3 and
# de-de!C3 - set every cell explicitly so all of them are updated.
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Error Detail (column P) messages for the failed handback.
$wsZhCn.Range("P3").Value = "Handback file name: rp0tgwtk.lxx is different with handoff file name: c2276190-ddf6-4c17-ba66-12f908df8e3e.2a9f85e9a9a88d85b23af682490470e3d654c446.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: rp0tgwtk.lxx is different with handoff file name: c2276190-ddf6-4c17-ba66-12f908df8e3e.2a9f85e9a9a88d85b23af682490470e3d654c446.de-de."

# Widen the Error Detail column so the new message is readable. The COM
# layer stores ColumnWidth with a constant +5/6 character padding baked into
# the serialized OOXML <col width>, so back the padding out here to land on
# exactly width="40" in the saved file.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
